$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1917.6923
$ws.Range("I43").Value = 1848.6428
$ws.Range("J43").Value = 1998.25
$ws.Range("K43").Value = 1848.6428
$ws.Range("L43").Value = 1998.25
$ws.Range("M43").Value = -1779.6428
$ws.Range("N43").Value = -2136.25
$ws.Range("H48").Value = 2880
$ws.Range("J48").Value = 2880
$ws.Range("L48").Value = 8640
$ws.Range("N48").Value = -9224
$ws.Range("H53").Value = 354.1
$ws.Range("I53").Value = 360.9
$ws.Range("J53").Value = 350.7
$ws.Range("K53").Value = 360.9
$ws.Range("L53").Value = 350.7
$ws.Range("M53").Value = 276.1
$ws.Range("N53").Value = -1624.7
$ws.Range("H56").Value = 2880
$ws.Range("J56").Value = 2880
$ws.Range("L56").Value = 8640
$ws.Range("N56").Value = -9708
$ws.Range("H81").Value = 38566.668
$ws.Range("J81").Value = 38566.668
$ws.Range("L81").Value = 38566.668
$ws.Range("N81").Value = -40562.668
$ws.Range("H84").Value = 38566.668
$ws.Range("J84").Value = 38566.668
$ws.Range("L84").Value = 115700.004
$ws.Range("N84").Value = -125684.004
$ws.Range("H134").Value = 61385
$ws.Range("J134").Value = 61385
$ws.Range("L134").Value = 61385
$ws.Range("N134").Value = -71525
$ws.Range("H137").Value = 1407.8379
$ws.Range("J137").Value = 2309.875
$ws.Range("L137").Value = 6929.625
$ws.Range("N137").Value = -12029.625
$ws.Range("H138").Value = 5609.965
$ws.Range("J138").Value = 6780.231
$ws.Range("L138").Value = 20340.693
$ws.Range("N138").Value = -30620.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 13969.923
$ws.Range("J6").Value = 11509.728
$ws.Range("L6").Value = 11509.728
$ws.Range("N6").Value = -11855.728
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4713
$ws.Range("N16").ClearContents()
$ws.Range("H26").Value = 3503.6155
$ws.Range("I26").Value = 2145.5833
$ws.Range("J26").Value = 19800
$ws.Range("K26").Value = 2145.5833
$ws.Range("L26").Value = 19800
$ws.Range("M26").Value = -1815.5833
$ws.Range("N26").Value = -20460
$ws.Range("H32").Value = 25198.727
$ws.Range("I32").Value = 4538.6885
$ws.Range("J32").Value = 130220.586
$ws.Range("K32").Value = 4538.6885
$ws.Range("L32").Value = 130220.586
$ws.Range("M32").Value = -4251.6885
$ws.Range("N32").Value = -130794.586
$ws.Range("H39").Value = 18400
$ws.Range("I39").Value = 7000
$ws.Range("J39").Value = 29800
$ws.Range("K39").Value = 7000
$ws.Range("L39").Value = 29800
$ws.Range("M39").Value = -6480
$ws.Range("N39").Value = -30840
$ws.Range("H61").Value = 1705.5245
$ws.Range("I61").Value = 1172.2565
$ws.Range("J61").Value = 2650.8635
$ws.Range("K61").Value = 1172.2565
$ws.Range("L61").Value = 2650.8635
$ws.Range("M61").Value = -960.2565
$ws.Range("N61").Value = -3074.8635
$ws.Range("H133").Value = 70000
$ws.Range("I133").Value = 70000
$ws.Range("K133").Value = 70000
$ws.Range("M133").Value = -67470
$ws.Range("H136").Value = 1705.5245
$ws.Range("I136").Value = 1172.2565
$ws.Range("J136").Value = 2650.8635
$ws.Range("K136").Value = 3516.7695
$ws.Range("L136").Value = 7952.5905
$ws.Range("M136").Value = -966.7694999999999
$ws.Range("N136").Value = -13052.5905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2990
$ws.Range("I15").Value = 1980
$ws.Range("J15").Value = 4000
$ws.Range("K15").Value = 1980
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = -1810
$ws.Range("N15").Value = -4340
$ws.Range("H31").Value = 26706.682
$ws.Range("I31").Value = 923.6539
$ws.Range("J31").Value = 44824.49
$ws.Range("K31").Value = 923.6539
$ws.Range("L31").Value = 44824.49
$ws.Range("M31").Value = -628.6539
$ws.Range("N31").Value = -45414.49
$ws.Range("H34").Value = 26706.682
$ws.Range("I34").Value = 923.6539
$ws.Range("J34").Value = 44824.49
$ws.Range("K34").Value = 923.6539
$ws.Range("L34").Value = 44824.49
$ws.Range("M34").Value = -721.6539
$ws.Range("N34").Value = -45228.49
$ws.Range("H50").Value = 9504.444
$ws.Range("I50").Value = 1000
$ws.Range("J50").Value = 10567.5
$ws.Range("K50").Value = 1000
$ws.Range("L50").Value = 10567.5
$ws.Range("M50").Value = -375
$ws.Range("N50").Value = -11817.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 326.25
$ws.Range("J11").Value = 430
$ws.Range("L11").Value = 1290
$ws.Range("N11").Value = -1570
$ws.Range("H75").Value = 2491.6667
$ws.Range("I75").Value = 950
$ws.Range("J75").Value = 2800
$ws.Range("K75").Value = 2850
$ws.Range("L75").Value = 8400
$ws.Range("M75").Value = -1852
$ws.Range("N75").Value = -10396
$ws.Range("H78").Value = 2491.6667
$ws.Range("I78").Value = 950
$ws.Range("J78").Value = 2800
$ws.Range("K78").Value = 8550
$ws.Range("L78").Value = 25200
$ws.Range("M78").Value = -3558
$ws.Range("N78").Value = -35184
$ws.Range("H122").Value = 779
$ws.Range("J122").Value = 948.5
$ws.Range("L122").Value = 8536.5
$ws.Range("N122").Value = -13436.5
$ws.Range("H141").Value = 3815.2727
$ws.Range("I141").Value = 1482.5
$ws.Range("J141").Value = 5148.2856
$ws.Range("K141").Value = 4447.5
$ws.Range("L141").Value = 15444.8568
$ws.Range("M141").Value = 732.5
$ws.Range("N141").Value = -25804.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1570.6571
$ws.Range("I132").Value = 1178.0358
$ws.Range("J132").Value = 3141.1428
$ws.Range("K132").Value = 3534.1074
$ws.Range("L132").Value = 9423.428400000001
$ws.Range("M132").Value = -1004.1074
$ws.Range("N132").Value = -14483.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2250
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 17161.334
$ws.Range("H141").Value = 57848.332
$ws.Range("J141").Value = 57848.332
$ws.Range("L141").Value = 57848.332
$ws.Range("N141").Value = -68208.33199999999
